$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 8752.75
$ws.Range("I21").Value = 8752.75
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 8752.75
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -8284.75
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 8752.75
$ws.Range("I23").Value = 8752.75
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 8752.75
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -8518.75
$ws.Range("N23").ClearContents()
$ws.Range("H33").Value = 5803.846
$ws.Range("I33").Value = 7308
$ws.Range("J33").Value = 3397.2
$ws.Range("K33").Value = 7308
$ws.Range("L33").Value = 3397.2
$ws.Range("M33").Value = -7079
$ws.Range("N33").Value = -3855.2
$ws.Range("H96").Value = 1676.2858
$ws.Range("I96").Value = 1141
$ws.Range("K96").Value = 3423
$ws.Range("M96").Value = -2050
$ws.Range("H98").Value = 2094.8235
$ws.Range("I98").Value = 786.7143
$ws.Range("J98").Value = 8199.333000000001
$ws.Range("K98").Value = 786.7143
$ws.Range("L98").Value = 8199.333000000001
$ws.Range("M98").Value = 711.2857
$ws.Range("N98").Value = -11195.333
$ws.Range("H106").Value = 6586
$ws.Range("I106").Value = 6913.5713
$ws.Range("K106").Value = 6913.5713
$ws.Range("M106").Value = -6282.5713
$ws.Range("H107").Value = 1818.8182
$ws.Range("I107").Value = 1478
$ws.Range("K107").Value = 1478
$ws.Range("M107").Value = 442
$ws.Range("H116").Value = 4851.5713
$ws.Range("I116").Value = 4928.357
$ws.Range("J116").Value = 4698
$ws.Range("K116").Value = 4928.357
$ws.Range("L116").Value = 4698
$ws.Range("M116").Value = -1486.357
$ws.Range("N116").Value = -11582
$ws.Range("H122").Value = 2094.8235
$ws.Range("I122").Value = 786.7143
$ws.Range("J122").Value = 8199.333000000001
$ws.Range("K122").Value = 2360.1429
$ws.Range("L122").Value = 24597.999
$ws.Range("M122").Value = 89.85710000000017
$ws.Range("N122").Value = -29497.999
$ws.Range("H132").Value = 3168.8647
$ws.Range("I132").Value = 3256.3142
$ws.Range("K132").Value = 9768.942599999998
$ws.Range("M132").Value = -7238.942599999998
$ws.Range("H135").Value = 1196.1666
$ws.Range("I135").Value = 814.11536
$ws.Range("K135").Value = 7327.03824
$ws.Range("M135").Value = -4792.03824
$ws.Range("H137").Value = 1689.4
$ws.Range("I137").Value = 1646.4117
$ws.Range("J137").Value = 1933
$ws.Range("K137").Value = 4939.2351
$ws.Range("L137").Value = 5799
$ws.Range("M137").Value = -2389.2351
$ws.Range("N137").Value = -10899
$ws.Range("H138").Value = 4532.0483
$ws.Range("I138").Value = 864.7059
$ws.Range("J138").Value = 8985.25
$ws.Range("K138").Value = 2594.1177
$ws.Range("L138").Value = 26955.75
$ws.Range("M138").Value = 2545.8823
$ws.Range("N138").Value = -37235.75
$ws.Range("H141").Value = 569.9167
$ws.Range("I141").Value = 569.9167
$ws.Range("K141").Value = 1709.7501
$ws.Range("M141").Value = 3470.2499

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1161.1111
$ws.Range("I5").Value = 1256.25
$ws.Range("K5").Value = 1256.25
$ws.Range("M5").Value = -1144.25
$ws.Range("H32").Value = 6521.855
$ws.Range("I32").Value = 484.9661
$ws.Range("K32").Value = 484.9661
$ws.Range("M32").Value = -197.9661
$ws.Range("H61").Value = 5215.83
$ws.Range("I61").Value = 3712.225
$ws.Range("K61").Value = 3712.225
$ws.Range("M61").Value = -3500.225
$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80540
$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -81872
$ws.Range("H74").Value = 2959.9583
$ws.Range("I74").Value = 2752.2
$ws.Range("K74").Value = 2752.2
$ws.Range("M74").Value = -1878.2
$ws.Range("H77").Value = 2959.9583
$ws.Range("I77").Value = 2752.2
$ws.Range("K77").Value = 13761
$ws.Range("M77").Value = -9393
$ws.Range("H88").Value = 1210.9286
$ws.Range("I88").Value = 873
$ws.Range("K88").Value = 873
$ws.Range("M88").Value = -467
$ws.Range("H91").Value = 1210.9286
$ws.Range("I91").Value = 873
$ws.Range("K91").Value = 873
$ws.Range("M91").Value = 531
$ws.Range("H102").Value = 9164
$ws.Range("I102").Value = 9655.866
$ws.Range("J102").Value = 5475
$ws.Range("K102").Value = 9655.866
$ws.Range("L102").Value = 5475
$ws.Range("M102").Value = -8033.866
$ws.Range("N102").Value = -8719
$ws.Range("H132").Value = 3362.4795
$ws.Range("I132").Value = 3073.5652
$ws.Range("K132").Value = 9220.695599999999
$ws.Range("M132").Value = -6690.695599999999
$ws.Range("H136").Value = 5215.83
$ws.Range("I136").Value = 3712.225
$ws.Range("K136").Value = 11136.675
$ws.Range("M136").Value = -8586.674999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1161.1111
$ws.Range("I4").Value = 1256.25
$ws.Range("K4").Value = 1256.25
$ws.Range("M4").Value = -1141.25
$ws.Range("H5").Value = 7259.4
$ws.Range("I5").Value = 3767.3333
$ws.Range("J5").Value = 12497.5
$ws.Range("K5").Value = 3767.3333
$ws.Range("L5").Value = 12497.5
$ws.Range("M5").Value = -3654.3333
$ws.Range("N5").Value = -12723.5
$ws.Range("H95").Value = 19235.6
$ws.Range("J95").Value = 19235.6
$ws.Range("L95").Value = 19235.6
$ws.Range("N95").Value = -24727.6
$ws.Range("H107").Value = 3737.5833
$ws.Range("I107").Value = 3594.4
$ws.Range("J107").Value = 3839.8572
$ws.Range("K107").Value = 3594.4
$ws.Range("L107").Value = 3839.8572
$ws.Range("M107").Value = -1674.4
$ws.Range("N107").Value = -7679.8572
$ws.Range("H134").Value = 3425.359
$ws.Range("I134").Value = 3183.5
$ws.Range("K134").Value = 9550.5
$ws.Range("M134").Value = -7015.5
$ws.Range("H138").Value = 88879.836
$ws.Range("J138").Value = 88879.836
$ws.Range("L138").Value = 88879.836
$ws.Range("N138").Value = -99159.836

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 17835.334
$ws.Range("I3").Value = 16507
$ws.Range("J3").Value = 18499.5
$ws.Range("K3").Value = 16507
$ws.Range("L3").Value = 18499.5
$ws.Range("M3").Value = -16394
$ws.Range("N3").Value = -18725.5
$ws.Range("H10").Value = 1835.3334
$ws.Range("I10").Value = 753.5
$ws.Range("J10").Value = 3999
$ws.Range("K10").Value = 753.5
$ws.Range("L10").Value = 3999
$ws.Range("M10").Value = -614.5
$ws.Range("N10").Value = -4277
$ws.Range("H31").Value = 6269.829
$ws.Range("I31").Value = 6694.7
$ws.Range("K31").Value = 6694.7
$ws.Range("M31").Value = -6399.7
$ws.Range("H33").Value = 4832.25
$ws.Range("I33").Value = 4832.25
$ws.Range("K33").Value = 4832.25
$ws.Range("M33").Value = -4453.25
$ws.Range("H34").Value = 6269.829
$ws.Range("I34").Value = 6694.7
$ws.Range("K34").Value = 6694.7
$ws.Range("M34").Value = -6492.7
$ws.Range("H58").Value = 1735
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H94").Value = 1432.3334
$ws.Range("I94").Value = 999
$ws.Range("K94").Value = 999
$ws.Range("M94").Value = -548
$ws.Range("H99").Value = 4131.609
$ws.Range("I99").Value = 3386.2942
$ws.Range("J99").Value = 6243.3335
$ws.Range("K99").Value = 3386.2942
$ws.Range("L99").Value = 6243.3335
$ws.Range("M99").Value = -1888.2942
$ws.Range("N99").Value = -9239.333500000001
$ws.Range("H107").Value = 968.7646999999999
$ws.Range("J107").Value = 993.3333
$ws.Range("L107").Value = 993.3333
$ws.Range("N107").Value = -4833.3333
$ws.Range("H122").Value = 302894.3
$ws.Range("I122").Value = 503758.16
$ws.Range("J122").Value = 1598.5
$ws.Range("K122").Value = 1511274.48
$ws.Range("L122").Value = 4795.5
$ws.Range("M122").Value = -1508824.48
$ws.Range("N122").Value = -9695.5
$ws.Range("H126").Value = 4131.609
$ws.Range("I126").Value = 3386.2942
$ws.Range("J126").Value = 6243.3335
$ws.Range("K126").Value = 10158.8826
$ws.Range("L126").Value = 18730.0005
$ws.Range("M126").Value = -7688.882599999999
$ws.Range("N126").Value = -23670.0005
$ws.Range("H132").Value = 1654.9546
$ws.Range("I132").Value = 1640.9445
$ws.Range("J132").Value = 1718
$ws.Range("K132").Value = 4922.833500000001
$ws.Range("L132").Value = 5154
$ws.Range("M132").Value = -2392.833500000001
$ws.Range("N132").Value = -10214
$ws.Range("H134").Value = 4366.357
$ws.Range("J134").Value = 11483.75
$ws.Range("L134").Value = 34451.25
$ws.Range("N134").Value = -39521.25
$ws.Range("H136").Value = 1735
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H141").Value = 67954.09
$ws.Range("I141").Value = 42000
$ws.Range("J141").Value = 70549.5
$ws.Range("K141").Value = 42000
$ws.Range("L141").Value = 70549.5
$ws.Range("M141").Value = -36820
$ws.Range("N141").Value = -80909.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 799.5
$ws.Range("I5").Value = 820
$ws.Range("K5").Value = 2460
$ws.Range("M5").Value = -2348
$ws.Range("H41").Value = 3190
$ws.Range("I41").Value = 2285
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 6855
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = -6517
$ws.Range("N41").Value = -15676
$ws.Range("H58").Value = 7905.6665
$ws.Range("J58").Value = 7905.6665
$ws.Range("L58").Value = 23716.9995
$ws.Range("N58").Value = -23972.9995
$ws.Range("H97").Value = 771.5
$ws.Range("J97").Value = 893.5
$ws.Range("L97").Value = 2680.5
$ws.Range("N97").Value = -3672.5
$ws.Range("H98").Value = 1544.6666
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1544.6666
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 4633.9998
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -7629.9998
$ws.Range("H99").Value = 1994.3334
$ws.Range("I99").Value = 788.7143
$ws.Range("K99").Value = 2366.1429
$ws.Range("M99").Value = -120.1428999999998
$ws.Range("H113").Value = 1729.4286
$ws.Range("J113").Value = 1789.3846
$ws.Range("L113").Value = 5368.1538
$ws.Range("N113").Value = -9708.1538
$ws.Range("H135").Value = 799.5
$ws.Range("I135").Value = 820
$ws.Range("K135").Value = 7380
$ws.Range("M135").Value = -4845

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 25000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 25000
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H43").Value = 1732.5
$ws.Range("I43").Value = 1732.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1732.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1581.5
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 16923
$ws.Range("I46").Value = 5999.875
$ws.Range("J46").Value = 34400
$ws.Range("K46").Value = 5999.875
$ws.Range("L46").Value = 34400
$ws.Range("M46").Value = -5843.875
$ws.Range("N46").Value = -34712
$ws.Range("H57").Value = 59856
$ws.Range("J57").Value = 59856
$ws.Range("L57").Value = 59856
$ws.Range("N57").Value = -61496
$ws.Range("H80").Value = 3735.4119
$ws.Range("I80").Value = 2885.9285
$ws.Range("K80").Value = 2885.9285
$ws.Range("M80").Value = -1887.9285
$ws.Range("H83").Value = 3735.4119
$ws.Range("I83").Value = 2885.9285
$ws.Range("K83").Value = 14429.6425
$ws.Range("M83").Value = -9437.6425
$ws.Range("H113").Value = 3750
$ws.Range("J113").Value = 3750
$ws.Range("L113").Value = 3750
$ws.Range("N113").Value = -8090
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 3781.3542
$ws.Range("I132").Value = 4146.1353
$ws.Range("J132").Value = 2554.3635
$ws.Range("K132").Value = 12438.4059
$ws.Range("L132").Value = 7663.0905
$ws.Range("M132").Value = -9908.4059
$ws.Range("N132").Value = -12723.0905

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26031.857
$ws.Range("I7").Value = 28425.295
$ws.Range("J7").Value = 15859.75
$ws.Range("K7").Value = 28425.295
$ws.Range("L7").Value = 15859.75
$ws.Range("M7").Value = -28313.295
$ws.Range("N7").Value = -16083.75
$ws.Range("H22").Value = 2517.1428
$ws.Range("I22").Value = 2276.75
$ws.Range("K22").Value = 2276.75
$ws.Range("M22").Value = -1981.75
$ws.Range("H23").Value = 4000000
$ws.Range("I23").Value = 4000000
$ws.Range("K23").Value = 4000000
$ws.Range("M23").Value = -3999770
$ws.Range("H27").Value = 2517.1428
$ws.Range("I27").Value = 2276.75
$ws.Range("K27").Value = 2276.75
$ws.Range("M27").Value = -2169.75
$ws.Range("H40").Value = 2711.0908
$ws.Range("I40").Value = 3070.2222
$ws.Range("J40").Value = 1095
$ws.Range("K40").Value = 3070.2222
$ws.Range("L40").Value = 1095
$ws.Range("M40").Value = -2934.2222
$ws.Range("N40").Value = -1367
$ws.Range("H46").Value = 4565.9546
$ws.Range("I46").Value = 4014.5715
$ws.Range("J46").Value = 5530.875
$ws.Range("K46").Value = 4014.5715
$ws.Range("L46").Value = 5530.875
$ws.Range("M46").Value = -3826.5715
$ws.Range("N46").Value = -5906.875
$ws.Range("H55").Value = 1116.3529
$ws.Range("I55").Value = 1417.375
$ws.Range("K55").Value = 1417.375
$ws.Range("M55").Value = -1244.375
$ws.Range("H61").Value = 1527.3334
$ws.Range("I61").Value = 1321.3
$ws.Range("K61").Value = 1321.3
$ws.Range("M61").Value = -1119.3
$ws.Range("H82").Value = 990.75
$ws.Range("I82").Value = 732.8333
$ws.Range("J82").Value = 1248.6666
$ws.Range("K82").Value = 732.8333
$ws.Range("L82").Value = 1248.6666
$ws.Range("M82").Value = -371.8333
$ws.Range("N82").Value = -1970.6666
$ws.Range("H85").Value = 990.75
$ws.Range("I85").Value = 732.8333
$ws.Range("J85").Value = 1248.6666
$ws.Range("K85").Value = 732.8333
$ws.Range("L85").Value = 1248.6666
$ws.Range("M85").Value = 515.1667
$ws.Range("N85").Value = -3744.6666
$ws.Range("H100").Value = 5849.9165
$ws.Range("J100").Value = 6611.1113
$ws.Range("L100").Value = 6611.1113
$ws.Range("N100").Value = -7693.1113
$ws.Range("H113").Value = 1527.3334
$ws.Range("I113").Value = 1321.3
$ws.Range("K113").Value = 1321.3
$ws.Range("M113").Value = 848.7
$ws.Range("H122").Value = 5744.263
$ws.Range("I122").Value = 5896.778
$ws.Range("K122").Value = 17690.334
$ws.Range("M122").Value = -15240.334
$ws.Range("H126").Value = 26031.857
$ws.Range("I126").Value = 28425.295
$ws.Range("J126").Value = 15859.75
$ws.Range("K126").Value = 85275.88499999999
$ws.Range("L126").Value = 47579.25
$ws.Range("M126").Value = -82805.88499999999
$ws.Range("N126").Value = -52519.25
$ws.Range("H132").Value = 2239.7693
$ws.Range("I132").Value = 1912.4839
$ws.Range("K132").Value = 5737.4517
$ws.Range("M132").Value = -3207.4517
$ws.Range("H136").Value = 1857.625
$ws.Range("I136").Value = 1743.5454
$ws.Range("K136").Value = 5230.6362
$ws.Range("M136").Value = -2680.6362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 25052.938
$ws.Range("I18").Value = 14142
$ws.Range("J18").Value = 33539.223
$ws.Range("K18").Value = 14142
$ws.Range("L18").Value = 33539.223
$ws.Range("M18").Value = -13969
$ws.Range("N18").Value = -33885.223
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H107").Value = 5459.4
$ws.Range("I107").Value = 5214.6924
$ws.Range("K107").Value = 15644.0772
$ws.Range("M107").Value = -13724.0772
$ws.Range("H113").Value = 3092.125
$ws.Range("I113").Value = 3092.125
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 9276.375
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -7106.375
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 5599.5
$ws.Range("I126").Value = 6899.3335
$ws.Range("K126").Value = 20698.0005
$ws.Range("M126").Value = -18228.0005
$ws.Range("H132").Value = 1856.791
$ws.Range("I132").Value = 1946.8246
$ws.Range("J132").Value = 1343.6
$ws.Range("K132").Value = 5840.4738
$ws.Range("L132").Value = 4030.8
$ws.Range("M132").Value = -3310.4738
$ws.Range("N132").Value = -9090.799999999999
$ws.Range("H136").Value = 2595.9395
$ws.Range("I136").Value = 1960.1305
$ws.Range("J136").Value = 4058.3
$ws.Range("K136").Value = 5880.3915
$ws.Range("L136").Value = 12174.9
$ws.Range("M136").Value = -3330.3915
$ws.Range("N136").Value = -17274.9

Write-Host "Applied all changes"